$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Price (D) and Volume(1h) (E) columns for rows 2-51
# NumberFormat is forced to Text ("@") before assignment so that
# numeric-looking strings (e.g. "49.00", "0.0345") are stored as exact
# text rather than being parsed into floating point numbers, then
# ClearFormats() removes the temporary format so the cell style stays
# identical to the original (no explicit style index).

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '66.712.83'
$c.ClearFormats()
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = '  -0.46%  '
$c.ClearFormats()
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '3.093.00'
$c.ClearFormats()
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = '  -0.83%  '
$c.ClearFormats()
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = '  -0.11%  '
$c.ClearFormats()
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '576.31'
$c.ClearFormats()
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = '  -0.71%  '
$c.ClearFormats()
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '172.88'
$c.ClearFormats()
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = '  -0.34%  '
$c.ClearFormats()
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = '  -0.11%  '
$c.ClearFormats()
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '3.091.83'
$c.ClearFormats()
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = '  -0.74%  '
$c.ClearFormats()
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.514'
$c.ClearFormats()
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = '  -1.53%  '
$c.ClearFormats()
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '6.38'
$c.ClearFormats()
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = '  -0.36%  '
$c.ClearFormats()
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = '  -2.24%  '
$c.ClearFormats()
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.473'
$c.ClearFormats()
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = '  -1.47%  '
$c.ClearFormats()
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = '  -3.78%  '
$c.ClearFormats()
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '35.71'
$c.ClearFormats()
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = '  -3.89%  '
$c.ClearFormats()
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = '  -0.63%  '
$c.ClearFormats()
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '3.607.45'
$c.ClearFormats()
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = '  -0.80%  '
$c.ClearFormats()
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '66.660.09'
$c.ClearFormats()
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = '  -0.54%  '
$c.ClearFormats()
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = '  +4.07%  '
$c.ClearFormats()
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '6.98'
$c.ClearFormats()
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = '  -2.01%  '
$c.ClearFormats()
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '3.091.34'
$c.ClearFormats()
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = '  -0.90%  '
$c.ClearFormats()
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '483.07'
$c.ClearFormats()
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = '  -1.57%  '
$c.ClearFormats()
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = '  -1.98%  '
$c.ClearFormats()
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = '  -2.46%  '
$c.ClearFormats()
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '83.33'
$c.ClearFormats()
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = '  -0.95%  '
$c.ClearFormats()
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '12.63'
$c.ClearFormats()
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = '  -4.74%  '
$c.ClearFormats()
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = '  -3.05%  '
$c.ClearFormats()
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '10.04'
$c.ClearFormats()
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = '  -3.68%  '
$c.ClearFormats()
$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = '  +0.09%  '
$c.ClearFormats()
$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = '  +0.98%  '
$c.ClearFormats()
$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = '  -4.01%  '
$c.ClearFormats()
$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = '  -3.64%  '
$c.ClearFormats()
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '27.90'
$c.ClearFormats()
$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = '  -2.33%  '
$c.ClearFormats()
$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = '  -2.61%  '
$c.ClearFormats()
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '0.0₃0932'
$c.ClearFormats()
$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = '  -1.42%  '
$c.ClearFormats()
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.ClearFormats()
$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = '  -0.12%  '
$c.ClearFormats()
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '48.11'
$c.ClearFormats()
$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = '  +1.73%  '
$c.ClearFormats()
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '5.58'
$c.ClearFormats()
$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = '  -4.98%  '
$c.ClearFormats()
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.941'
$c.ClearFormats()
$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = '  -3.43%  '
$c.ClearFormats()
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '49.00'
$c.ClearFormats()
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = '  -2.32%  '
$c.ClearFormats()
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = '  -0.86%  '
$c.ClearFormats()
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = '  -0.73%  '
$c.ClearFormats()
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '1.96'
$c.ClearFormats()
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = '  -4.67%  '
$c.ClearFormats()
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = '  -3.00%  '
$c.ClearFormats()
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '2.62'
$c.ClearFormats()
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = '  +0.60%  '
$c.ClearFormats()
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '2.795.48'
$c.ClearFormats()
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = '  -1.19%  '
$c.ClearFormats()
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.0345'
$c.ClearFormats()
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '368.10'
$c.ClearFormats()
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = '  -4.00%  '
$c.ClearFormats()
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = '  -0.86%  '
$c.ClearFormats()
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = '  -0.01%  '
$c.ClearFormats()
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '24.41'
$c.ClearFormats()
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = '  -2.34%  '
$c.ClearFormats()
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '2.18'
$c.ClearFormats()
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = '  -1.23%  '
$c.ClearFormats()
